$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply crypto price/volume updates per commit diff.
# D-column cells use a leading apostrophe to force text storage
# (matching the original inlineStr text cells), avoiding Excel
# auto-converting numeric-looking strings into numbers; the Style
# is reset to Normal afterward so no stray formatting is introduced.

$ws.Range("D2").Value = "'27.129.76"

$ws.Range("D3").Value = "'1.568.23"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("D5").Value = "'210.78"
$ws.Range("E5").Value = "  +2.03%  "

$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("D10").Value = "'0.0598"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "'1.788.93"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").Value = "'1.523.75"
$ws.Range("E13").Value = "  -1.75%  "

$ws.Range("D14").Value = "'3.78"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").Value = "'27.120.23"

$ws.Range("D17").Value = "'62.05"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "'0.0₃0704"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'214.86"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.38"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").Value = "'9.21"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "'154.20"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("D27").Value = "'15.06"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("E30").Value = "  +5.12%  "

$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").Value = "  +3.00%  "

$ws.Range("D34").Value = "'1.430.38"
$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("E35").Value = "  +10.75%  "

$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "'5.85"
$ws.Range("E40").Value = "  +3.46%  "

$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("E43").Value = "  +2.49%  "

$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").Value = "'64.62"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  +0.94%  "

$ws.Range("D47").Value = "'1.706.97"

$ws.Range("D48").Value = "'86.07"
$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0102"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0518"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "  +0.14%  "

# Reset style on all touched cells to avoid leftover formatting artifacts.
foreach ($addr in @("D2","D3","E3","D5","E5","E7","E8","E9","D10","E10","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","D17","E17","D18","E18","B19","C19","D19","E19","B20","C20","D20","E20","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","E28","E29","E30","E31","E32","E33","D34","E34","E35","E36","E37","E38","E39","D40","E40","E41","E42","E43","E44","D45","E45","E46","D47","D48","E48","B49","C49","D49","E49","B50","C50","D50","E50","D51","E51")) {
    $ws.Range($addr).Style = "Normal"
}
